# Auto-generated edit script applying diff changes to Sargatanas_Profits workbook
# Sheet tabs: ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H3").Value = 52999
$ws.Range("J3").Value = 52999
$ws.Range("L3").Value = 52999
$ws.Range("N3").Value = -53227
$ws.Range("H4").Value = 710.2727
$ws.Range("I4").Value = 749.3333
$ws.Range("K4").Value = 749.3333
$ws.Range("M4").Value = -635.3333
$ws.Range("H33").Value = 463.44446
$ws.Range("I33").Value = 95.69231000000001
$ws.Range("J33").Value = 1419.6
$ws.Range("K33").Value = 95.69231000000001
$ws.Range("L33").Value = 1419.6
$ws.Range("M33").Value = 133.30769
$ws.Range("N33").Value = -1877.6
$ws.Range("H43").Value = 686158.2
$ws.Range("J43").Value = 1027349.75
$ws.Range("L43").Value = 1027349.75
$ws.Range("N43").Value = -1027487.75
$ws.Range("H102").Value = 52999
$ws.Range("J102").Value = 52999
$ws.Range("L102").Value = 52999
$ws.Range("N102").Value = -59489
$ws.Range("H138").Value = 8255.286
$ws.Range("I138").Value = 4306.636
$ws.Range("J138").Value = 10810.294
$ws.Range("K138").Value = 12919.908
$ws.Range("L138").Value = 32430.882
$ws.Range("M138").Value = -7779.908000000001
$ws.Range("N138").Value = -42710.882

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2072618.8
$ws.Range("I32").Value = 2072618.8
$ws.Range("K32").Value = 2072618.8
$ws.Range("M32").Value = -2072331.8
$ws.Range("H45").Value = 2910.5
$ws.Range("I45").Value = 1426.8334
$ws.Range("K45").Value = 1426.8334
$ws.Range("M45").Value = -1049.8334
$ws.Range("H61").Value = 62512156
$ws.Range("I61").Value = 10500
$ws.Range("K61").Value = 10500
$ws.Range("M61").Value = -10288
$ws.Range("H132").Value = 7347.2188
$ws.Range("I132").Value = 4252.5835
$ws.Range("K132").Value = 12757.7505
$ws.Range("M132").Value = -10227.7505
$ws.Range("H136").Value = 62512156
$ws.Range("I136").Value = 10500
$ws.Range("K136").Value = 31500
$ws.Range("M136").Value = -28950

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 6415813.5
$ws.Range("I20").Value = 8775887
$ws.Range("K20").Value = 8775887
$ws.Range("M20").Value = -8775640
$ws.Range("H29").Value = 3666.6667
$ws.Range("I29").Value = 3750
$ws.Range("J29").Value = 3500
$ws.Range("K29").Value = 3750
$ws.Range("L29").Value = 3500
$ws.Range("M29").Value = -3461
$ws.Range("N29").Value = -4078
$ws.Range("H39").Value = 0
$ws.Range("J39").Value = 0
$ws.Range("L39").Value = 0
$ws.Range("N39").ClearContents()
$ws.Range("H111").Value = 80222.5
$ws.Range("J111").Value = 80222.5
$ws.Range("L111").Value = 80222.5
$ws.Range("N111").Value = -88402.5
$ws.Range("H123").Value = 76997.5
$ws.Range("J123").Value = 76997.5
$ws.Range("L123").Value = 76997.5
$ws.Range("N123").Value = -86797.5
$ws.Range("H134").Value = 8072115.5
$ws.Range("I134").Value = 20836876
$ws.Range("K134").Value = 62510628
$ws.Range("M134").Value = -62508093

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H6").Value = 11833.333
$ws.Range("J6").Value = 9000
$ws.Range("L6").Value = 9000
$ws.Range("N6").Value = -9226
$ws.Range("H31").Value = 7691.909
$ws.Range("I31").Value = 2902
$ws.Range("J31").Value = 11683.5
$ws.Range("K31").Value = 2902
$ws.Range("L31").Value = 11683.5
$ws.Range("M31").Value = -2607
$ws.Range("N31").Value = -12273.5
$ws.Range("H34").Value = 7691.909
$ws.Range("I34").Value = 2902
$ws.Range("J34").Value = 11683.5
$ws.Range("K34").Value = 2902
$ws.Range("L34").Value = 11683.5
$ws.Range("M34").Value = -2700
$ws.Range("N34").Value = -12087.5
$ws.Range("H43").Value = 0
$ws.Range("J43").Value = 0
$ws.Range("L43").Value = 0
$ws.Range("N43").ClearContents()
$ws.Range("H58").Value = 5857.6924
$ws.Range("I58").Value = 2890.5312
$ws.Range("J58").Value = 10605.15
$ws.Range("K58").Value = 2890.5312
$ws.Range("L58").Value = 10605.15
$ws.Range("M58").Value = -2687.5312
$ws.Range("N58").Value = -11011.15
$ws.Range("H80").Value = 52000
$ws.Range("J80").Value = 52000
$ws.Range("L80").Value = 52000
$ws.Range("N80").Value = -54246
$ws.Range("H83").Value = 52000
$ws.Range("J83").Value = 52000
$ws.Range("L83").Value = 156000
$ws.Range("N83").Value = -167232
$ws.Range("H101").Value = 0
$ws.Range("J101").Value = 0
$ws.Range("L101").Value = 0
$ws.Range("N101").ClearContents()
$ws.Range("H132").Value = 7539.875
$ws.Range("I132").Value = 5295.5264
$ws.Range("J132").Value = 10820.077
$ws.Range("K132").Value = 15886.5792
$ws.Range("L132").Value = 32460.231
$ws.Range("M132").Value = -13356.5792
$ws.Range("N132").Value = -37520.231
$ws.Range("H134").Value = 4614.86
$ws.Range("I134").Value = 1783.5
$ws.Range("K134").Value = 5350.5
$ws.Range("M134").Value = -2815.5
$ws.Range("H136").Value = 5857.6924
$ws.Range("I136").Value = 2890.5312
$ws.Range("J136").Value = 10605.15
$ws.Range("K136").Value = 8671.5936
$ws.Range("L136").Value = 31815.45
$ws.Range("M136").Value = -6121.5936
$ws.Range("N136").Value = -36915.45
$ws.Range("H141").Value = 68662
$ws.Range("J141").Value = 70613.71000000001
$ws.Range("L141").Value = 70613.71000000001
$ws.Range("N141").Value = -80973.71000000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 77580.234
$ws.Range("I2").Value = 69.416664
$ws.Range("J2").Value = 144018.08
$ws.Range("K2").Value = 416.499984
$ws.Range("L2").Value = 864108.48
$ws.Range("M2").Value = -303.499984
$ws.Range("N2").Value = -864334.48
$ws.Range("H60").Value = 1260
$ws.Range("I60").Value = 500
$ws.Range("J60").Value = 1450
$ws.Range("K60").Value = 1500
$ws.Range("L60").Value = 4350
$ws.Range("M60").Value = -1249
$ws.Range("N60").Value = -4852
$ws.Range("H102").Value = 11400
$ws.Range("I102").Value = 11400
$ws.Range("K102").Value = 34200
$ws.Range("M102").Value = -31766
$ws.Range("H129").Value = 134491.33
$ws.Range("J129").Value = 500499.5
$ws.Range("L129").Value = 1501498.5
$ws.Range("N129").Value = -1511498.5
$ws.Range("H131").Value = 78862.16
$ws.Range("J131").Value = 335163.66
$ws.Range("L131").Value = 1005490.98
$ws.Range("N131").Value = -1015570.98

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H15").Value = 40000
$ws.Range("J15").Value = 40000
$ws.Range("L15").Value = 40000
$ws.Range("N15").Value = -40576
$ws.Range("H81").Value = 40000
$ws.Range("J81").Value = 40000
$ws.Range("L81").Value = 40000
$ws.Range("N81").Value = -41996
$ws.Range("H84").Value = 40000
$ws.Range("J84").Value = 40000
$ws.Range("L84").Value = 120000
$ws.Range("N84").Value = -129984
$ws.Range("H112").Value = 0
$ws.Range("J112").Value = 0
$ws.Range("L112").Value = 0
$ws.Range("N112").ClearContents()
$ws.Range("H125").Value = 85000
$ws.Range("J125").Value = 85000
$ws.Range("L125").Value = 85000
$ws.Range("N125").Value = -89920
$ws.Range("H132").Value = 6684.8623
$ws.Range("I132").Value = 4450.7393
$ws.Range("K132").Value = 13352.2179
$ws.Range("M132").Value = -10822.2179

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 7199.8
$ws.Range("I7").Value = 6666.3335
$ws.Range("K7").Value = 6666.3335
$ws.Range("M7").Value = -6554.3335
$ws.Range("H55").Value = 705.26666
$ws.Range("I55").Value = 531.1667
$ws.Range("J55").Value = 821.3333
$ws.Range("K55").Value = 531.1667
$ws.Range("L55").Value = 821.3333
$ws.Range("M55").Value = -358.1667
$ws.Range("N55").Value = -1167.3333
$ws.Range("H122").Value = 5050.4194
$ws.Range("I122").Value = 4395.9653
$ws.Range("K122").Value = 13187.8959
$ws.Range("M122").Value = -10737.8959
$ws.Range("H126").Value = 7199.8
$ws.Range("I126").Value = 6666.3335
$ws.Range("K126").Value = 19999.0005
$ws.Range("M126").Value = -17529.0005
$ws.Range("H132").Value = 9440732
$ws.Range("I132").Value = 12199482
$ws.Range("K132").Value = 36598446
$ws.Range("M132").Value = -36595916
$ws.Range("H136").Value = 13732.467
$ws.Range("I136").Value = 6289.9165
$ws.Range("J136").Value = 22238.238
$ws.Range("K136").Value = 18869.7495
$ws.Range("L136").Value = 66714.71400000001
$ws.Range("M136").Value = -16319.7495
$ws.Range("N136").Value = -71814.71400000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 24285.666
$ws.Range("I113").Value = 44932.5
$ws.Range("J113").Value = 3638.8333
$ws.Range("K113").Value = 134797.5
$ws.Range("L113").Value = 10916.4999
$ws.Range("M113").Value = -132627.5
$ws.Range("N113").Value = -15256.4999
$ws.Range("H132").Value = 9029.875
$ws.Range("I132").Value = 11849.863
$ws.Range("J132").Value = 5583.222
$ws.Range("K132").Value = 35549.589
$ws.Range("L132").Value = 16749.666
$ws.Range("M132").Value = -33019.589
$ws.Range("N132").Value = -21809.666
$ws.Range("H136").Value = 17725454
$ws.Range("I136").Value = 29416932
$ws.Range("K136").Value = 88250796
$ws.Range("M136").Value = -88248246
